$d = $word.ActiveDocument

# 1. Remove the stray _GoBack bookmark that sat at the end of the "Agenda" heading.
$d.Bookmarks.Item("_GoBack").Delete()

# 2. Mention the new laser distance measurement IC component alongside the embedded
#    arm camera bullet.
$d.Content.Find.Execute( `
    "Embedded arm camera, motors and", $true, $false, $false, $false, $false, `
    $true, 1, $false, `
    "Embedded arm camera, laser distance measurement IC, motors and", 2)

# 3. Mention the same component in the detachable end-effector power description.
$d.Content.Find.Execute( `
    "to it to power small motors, camera, LEDs and maybe a small microcontroller.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "to it to power small motors, laser distance measurement IC, camera, LEDs and maybe a small microcontroller.", 2)

# 4. The _GoBack bookmark re-appears mid-word in the "feasibility and" bullet,
#    reflecting where the author's cursor was left after their last edit.
$rng = $d.Content
$rng.Find.Execute("Look into the feasibility a", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng)

# 5. Add three new bullets about a temperature-controlled fan switch after the
#    "cables to be correct" bullet.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*This will allow the design for cables to be correct*") {
        $target = $p
    }
}

$target.Range.InsertParagraphAfter()
$p1 = $target.Next()
$p1.Range.Text = "Create a temp-dependant switch to turn on the fans"
$p1.Range.ListFormat.ListLevelNumber = 2

$p1.Range.InsertParagraphAfter()
$p2 = $p1.Next()
$p2.Range.Text = "Research into an appropriate sensor"
$p2.Range.ListFormat.ListLevelNumber = 3

$p2.Range.InsertParagraphAfter()
$p3 = $p2.Next()
$p3.Range.Text = "Wire up something to turn the fans on when a certain temperature occurs"
$p3.Range.ListFormat.ListLevelNumber = 3
